$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wood Pellets row (row 9): unit changes from "ton" to "pound",
#     BTUs per unit changes from 16,000,000 to 8,300, and CO2 changes
#     from 67 to 62.9 ---
$ws.Range("C9").Value = "pound"
$ws.Range("D9").Value = 8300
$ws.Range("E9").Value = 62.9

# --- New comments describing the btus (D) and co2 (E) header columns ---
$ws.Range("D1").AddComment("BTUs per Fuel Unit (e.g. gallon)")
$ws.Range("E1").AddComment("Pounds / MMBTU")

# --- New comment explaining the new Wood Pellets BTU figure ---
$ws.Range("D9").AddComment("ChatGPT says 8000 – 8600 BTU / lb")

# --- Updated comment explaining the new Wood Pellets CO2 figure ---
[void]$ws.Range("E9").Comment.Text("7% moisture content (based on dry weight). 0.508 lb Carbon per lb dry wood. 70% decay offset.")

# --- Restore selection/view state to reflect the last edited cell ---
[void]$ws.Range("H13").Select()
